# Applies the OOXML diff:
#  1. "O" + "ffentlig IP-adresse" -> merged run "Offentlig IP-adresse"
#  2. "søgning" run gains w:lang val="en-US"
#  3. CMD nslookup runs merged into a single run (drops spell-check markers)
#  4. CMD ping studypoints.info runs merged into a single run
#  5. Applikationslaget paragraph: several runs merged together
#  6. HTTP / Hypertext runs merged into a single run
#  7. TCP paragraph: final answer replaced with new text

$d = $word.ActiveDocument

# 1. Merge "O" / "ffentlig IP-adresse" into "Offentlig IP-adresse"
$d.Content.Find.Execute("Offentlig IP-adresse", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Offentlig IP-adresse", 2) | Out-Null

# 2. Set language of the "søgning" run to en-US
$r = $d.Content
$r.Find.Execute("søgning", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.LanguageID = "en-US"

# 3. Merge CMD: “nslookup myip.opendns.com resolver1.opendns.com” into one run
$d.Content.Find.Execute("CMD: “nslookup myip.opendns.com resolver1.opendns.com”", $true, $false, $false, $false, $false, `
    $true, 1, $false, "CMD: “nslookup myip.opendns.com resolver1.opendns.com”", 2) | Out-Null

# 4. Merge CMD: “ping studypoints.info” into one run
$d.Content.Find.Execute("CMD: “ping studypoints.info”", $true, $false, $false, $false, $false, `
    $true, 1, $false, "CMD: “ping studypoints.info”", 2) | Out-Null

# 5. Merge the long Applikationslaget explanation into fewer runs
$d.Content.Find.Execute("Applikationslaget samler data fra software og pakker det løst samlet. Transportlaget samler det data som den har modtaget fra applikationslaget og smækker det i en “pakke”. I netværkslaget bliver tilføjet ip-adresse til “pakken”. I linklaget f", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Applikationslaget samler data fra software og pakker det løst samlet. Transportlaget samler det data som den har modtaget fra applikationslaget og smækker det i en “pakke”. I netværkslaget bliver tilføjet ip-adresse til “pakken”. I linklaget f", 2) | Out-Null

$d.Content.Find.Execute("r “pakke", $true, $false, $false, $false, $false, `
    $true, 1, $false, "r “pakke", 2) | Out-Null

$d.Content.Find.Execute("” at vide hvilken mac-adresse den skal til. Sidst er det ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "” at vide hvilken mac-adresse den skal til. Sidst er det ", 2) | Out-Null

$d.Content.Find.Execute(" som man ikke altid tæller med, men det er i bund og grund den fysiske forbindelse mellem senderen og modtageren", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " som man ikke altid tæller med, men det er i bund og grund den fysiske forbindelse mellem senderen og modtageren", 2) | Out-Null

# 6. Merge "HTTP står for “" + "Hypertext" + " Transfer Protocol..." into one run
$d.Content.Find.Execute("HTTP står for “Hypertext Transfer Protocol” og bruges til at sende og modtage datapakker i forbindelse med almindeligt webbrug.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "HTTP står for “Hypertext Transfer Protocol” og bruges til at sende og modtage datapakker i forbindelse med almindeligt webbrug.", 2) | Out-Null

# 7. Replace the final TCP answer with the new text
$d.Content.Find.Execute("TCP holder forbindelsen aktiv mellem serveren og klienten, og derved kan vi eventuelt oprette sessions i vores webapplikationer, hvis dette er nødvendigt.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "HTTP er en del af applikationslaget, og TCP er en del af “transportlaget” hvilket vil sige, at TCP står for at sende HTTP datapakkerne.", 2) | Out-Null

Write-Host "done"
